# Fruta / hortaliza, semanal
# Insert two new weekly price records (rows 270-271) into the
# "Feria Lagunitas de Puerto Montt - Ciruela" sheet, pushing the
# existing rows 270-305 down to 272-307.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 270 (shifts everything below down by 2)
$ws.Rows.Item(270).Insert()
$ws.Rows.Item(270).Insert()

# New row 270: Angeleno / Primera
$ws.Range("A270").Value2 = 4
$ws.Range("B270").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C270").Value2 = "Los Lagos"
$ws.Range("D270").Value2 = 44995
$ws.Range("E270").Value2 = 10
$ws.Range("F270").Value2 = "Fruta"
$ws.Range("G270").Value2 = 100103
$ws.Range("H270").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I270").Value2 = 100103002
$ws.Range("J270").Value2 = "Ciruela"
$ws.Range("K270").Value2 = "Angeleno"
$ws.Range("L270").Value2 = "Primera"
$ws.Range("M270").Value2 = 600
$ws.Range("N270").Value2 = 14000
$ws.Range("O270").Value2 = 15000
$ws.Range("P270").Value2 = 14500
$ws.Range("Q270").Value2 = "$/caja 15 kilos granel"
$ws.Range("R270").Value2 = "Región de O'Higgins"
$ws.Range("S270").Value2 = 1036
$ws.Range("T270").Value2 = 14

# New row 271: Angeleno / Segunda
$ws.Range("A271").Value2 = 4
$ws.Range("B271").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C271").Value2 = "Los Lagos"
$ws.Range("D271").Value2 = 44995
$ws.Range("E271").Value2 = 10
$ws.Range("F271").Value2 = "Fruta"
$ws.Range("G271").Value2 = 100103
$ws.Range("H271").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I271").Value2 = 100103002
$ws.Range("J271").Value2 = "Ciruela"
$ws.Range("K271").Value2 = "Angeleno"
$ws.Range("L271").Value2 = "Segunda"
$ws.Range("M271").Value2 = 300
$ws.Range("N271").Value2 = 11000
$ws.Range("O271").Value2 = 11000
$ws.Range("P271").Value2 = 11000
$ws.Range("Q271").Value2 = "$/caja 14 kilos granel"
$ws.Range("R271").Value2 = "Región de O'Higgins"
$ws.Range("S271").Value2 = 786
$ws.Range("T271").Value2 = 14
